$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.213.75"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").Value = "1.860.64"
$ws.Range("E3").Value = "  -0.60%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9993"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7140"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "240.44"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.59%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9992"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07751"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3090"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.13"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08287"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.54%  "
$ws.Range("D12").Value = "1.885.43"
$ws.Range("E12").Value = "  +0.75%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.224"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7168"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.41%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "90.84"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.00%  "
$ws.Range("D16").Value = "29.198.78"
$ws.Range("E16").Value = "  -0.84%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.883"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "244.54"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007814"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.91%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.16"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.80%  "
$ws.Range("D21").Value = "2.109.28"
$ws.Range("E21").Value = "  -1.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9991"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.017"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +3.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9995"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1601"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +3.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.40"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.921"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.61"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.81%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.329"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.00%  "
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.462"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.82%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.497"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.243"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05183"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.70%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.917"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.170"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7324"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +2.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.666"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01851"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.693"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.83%  "
$ws.Range("D40").Value = "1.168.40"
$ws.Range("E40").Value = "  -2.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9055"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.146"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "72.86"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.81%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9986"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.83"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.47%  "
$ws.Range("D46").Value = "2.004.25"
$ws.Range("E46").Value = "  -1.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5214"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.771"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.02%  "
$ws.Range("E49").Value = "  +1.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.318"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.874"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.57%  "

Write-Output "done"
